$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("A1").Value = "Fund Name"
$ws.Range("B1").Value = "Rate of Return (per month)"
$ws.Range("C1").Value = "Time (in months)"
$ws.Range("D1").Value = "SIP "
$ws.Range("E1").Value = "Principle Investment"
$ws.Range("F1").Value = "Final Amount"

# Bold header row
$ws.Range("A1:F1").Font.Bold = $true

# Row 2 - X Fund
$ws.Range("A2").Value = "X Fund "
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = -1
$ws.Range("E2").Value = -100
$ws.Range("F2").Formula = "=FV(B2,C2,D2,E2)"

# Row 3 - Y_Fund
$ws.Range("A3").Value = "Y_Fund"
$ws.Range("B3").Value = 0.02
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = -3
$ws.Range("E3").Value = -50
$ws.Range("F3").Formula = "=FV(B3,C3,D3,E3)"

# Row 4 - Z_fund
$ws.Range("A4").Value = "Z_fund"
$ws.Range("B4").Value = 0.05
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = -1100
$ws.Range("E4").Value = 0
$ws.Range("F4").Formula = "=FV(B4,C4,D4,E4)"

# Number formats
$ws.Range("B2:B4").NumberFormat = "0%"
$ws.Range("F2:F4").NumberFormat = """₹"" #,##0.00;[Red]""₹"" -#,##0.00"

# Column widths
$ws.Columns.Item(2).ColumnWidth = 27.44140625
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 6
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(6).ColumnWidth = 11.88671875

# Selection to match diff (activeCell B4)
$ws.Range("B4").Select()

$ws.PageSetup.Orientation = 1
